$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Arkham Horror The Card Game: Revised Core Set"
$ws.Range("C17").Value = "https://www.amazon.de/-/en/Fantasy-Flight-Games-Arkham-Horror/dp/B0999D3P8S/ref=sr_1_3?crid=HGFOXKZQTV6Y&keywords=arkham+horror+english&qid=1699364257&sprefix=arkham+horror+en%2Caps%2C89&sr=8-3"
$ws.Range("D17").Value = "87.48 EUR"
$ws.Range("B17").Value = "https://m.media-amazon.com/images/I/61wUiOvhH2L._AC_SL1075_.jpg"

$ws.Range("B17").Select()
